$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 139, shifting existing rows 139-185 down to 140-186.
$ws.Rows.Item(139).Insert()

# Populate the new row 139 with the new data record.
$ws.Range("A139").Value = 8
$ws.Range("B139").Value = "Terminal La Palmera de La Serena"
$ws.Range("C139").Value = "Coquimbo"
$ws.Range("D139").Value = 44468
$ws.Range("E139").Value = 4
$ws.Range("F139").Value = 100112032
$ws.Range("G139").Value = "Zapallo italiano"
$ws.Range("H139").Value = "Sin especificar"
$ws.Range("I139").Value = "Primera"
$ws.Range("J139").Value = 500
$ws.Range("K139").Value = 16000
$ws.Range("L139").Value = 17000
$ws.Range("M139").Value = 16500
$ws.Range("N139").Value = "$/caja 70 unidades"
$ws.Range("O139").Value = "Provincia de Limarí"
$ws.Range("P139").Value = 236
$ws.Range("Q139").Value = 70
$ws.Range("R139").Value = "Hortaliza"
